# Fruta / hortaliza, semanal
# Insert a new data row at row 93 (shifting existing rows 93..175 down to 94..176)
# and populate it with the new daily price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 93; this shifts rows 93-175 down to 94-176
# and also updates the sheet dimension automatically (A1:T175 -> A1:T176).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new record's data.
$ws.Range("A93").Value = 11
$ws.Range("B93").Value = "Vega Monumental Concepción"
$ws.Range("C93").Value = "Bíobío"
$ws.Range("D93").Value = 44960
$ws.Range("E93").Value = 8
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100109
$ws.Range("H93").Value = "Uva"
$ws.Range("I93").Value = 100109001
$ws.Range("J93").Value = "Uva"
$ws.Range("K93").Value = "Flame Seedless"
$ws.Range("L93").Value = "Primera"
$ws.Range("M93").Value = 100
$ws.Range("N93").Value = 10000
$ws.Range("O93").Value = 11000
$ws.Range("P93").Value = 10500
$ws.Range("Q93").Value = "`$/bandeja 18 kilos"
$ws.Range("R93").Value = "Provincia de Limarí"
$ws.Range("S93").Value = 583
$ws.Range("T93").Value = 18
